$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

# Header rename
$ws.Range("B1").Value = "Non-residential"

# Updated projection values (new parquet files w/ 27 countries)
$ws.Range("B2").Value = 420.5118763716624
$ws.Range("C2").Value = 723.2113187875967

$ws.Range("B3").Value = 417.2098904849987
$ws.Range("C3").Value = 708.340439546219

$ws.Range("B4").Value = 415.72087345171656
$ws.Range("C4").Value = 699.7582047619915

$ws.Range("B5").Value = 414.8188235513056
$ws.Range("C5").Value = 691.7742393846

$ws.Range("B6").Value = 414.5569372331537
$ws.Range("C6").Value = 683.1864284357799

$ws.Range("B7").Value = 415.43283424500953
$ws.Range("C7").Value = 672.3298372900899

$ws.Range("B8").Value = 417.0907195523708
$ws.Range("C8").Value = 660.3146200847184
